# Finish processing the questions sheet: the "答案" (answer) column is
# redundant bookkeeping now that each option row already carries its own
# 选项序号 (option #), 选项 (option text) and 是否正确答案 (is-correct) data.
# Deleting the whole column shifts 选项序号/选项/是否正确答案 one slot to the
# left, which is exactly what the sheet should look like going forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Delete()

# Leave the view scrolled over to, and focused on, the option-text column
# the edits were made in.
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E7").Select()
